$wb = $excel.ActiveWorkbook

# Overview sheet - "Latest HO Xliff Generate Date" for the first data row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-30 03:09:27"

# zh-cn sheet - "Correspond Handoff Datetime" and "Correspond Handback DateTime" for the first data row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-30 03:09:22"
$wsZhCn.Range("K2").Value = "2016-08-30 03:09:40"

# de-de sheet - "Correspond Handback DateTime" for the first data row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-30 03:09:46"
